$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Update Constraint(s) for SubjectPerson.typeId (row 5, column AJ)
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}
"
